$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values
# for rows 2-51. Price cells that would otherwise be auto-parsed as numbers
# (losing trailing zeros, e.g. "1.00" -> 1) are written with a leading
# apostrophe via .Formula so Excel stores them as quote-prefixed text,
# matching the original inlineStr/text representation.

$ws.Range("D2").Value2 = "59.259.99"
$ws.Range("E2").Value2 = "  +5.31%  "
$ws.Range("D3").Value2 = "3.327.00"
$ws.Range("E3").Value2 = "  +3.27%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value2 = "  -0.20%  "
$ws.Range("D5").Formula = "'403.17"
$ws.Range("E5").Value2 = "  +1.28%  "
$ws.Range("D6").Formula = "'110.41"
$ws.Range("E6").Value2 = "  -0.50%  "
$ws.Range("D7").Formula = "'0.589"
$ws.Range("E7").Value2 = "  +6.26%  "
$ws.Range("E8").Value2 = "  -0.09%  "
$ws.Range("D9").Formula = "'0.637"
$ws.Range("E9").Value2 = "  +2.94%  "
$ws.Range("D10").Formula = "'39.86"
$ws.Range("E10").Value2 = "  +1.59%  "
$ws.Range("D11").Formula = "'0.0997"
$ws.Range("E11").Value2 = "  +7.83%  "
$ws.Range("E12").Value2 = "  +1.59%  "
$ws.Range("D13").Value2 = "3.833.47"
$ws.Range("E13").Value2 = "  +2.66%  "
$ws.Range("D14").Formula = "'8.42"
$ws.Range("E14").Value2 = "  +4.00%  "
$ws.Range("D15").Formula = "'19.38"
$ws.Range("E15").Value2 = "  +1.86%  "
$ws.Range("D16").Value2 = "3.320.24"
$ws.Range("E16").Value2 = "  +2.89%  "
$ws.Range("E17").Value2 = "  +0.14%  "
$ws.Range("D18").Value2 = "58.805.98"
$ws.Range("E18").Value2 = "  +4.64%  "
$ws.Range("D19").Formula = "'10.89"
$ws.Range("E19").Value2 = "  -0.58%  "
$ws.Range("E20").Value2 = "  -0.82%  "
$ws.Range("D21").Formula = "'0.0000111"
$ws.Range("E21").Value2 = "  +6.76%  "
$ws.Range("D22").Formula = "'13.06"
$ws.Range("E22").Value2 = "  +0.24%  "
$ws.Range("D23").Formula = "'305.80"
$ws.Range("E23").Value2 = "  +2.97%  "
$ws.Range("D24").Formula = "'75.07"
$ws.Range("E24").Value2 = "  -1.04%  "
$ws.Range("D25").Formula = "'3.20"
$ws.Range("E25").Value2 = "  -0.52%  "
$ws.Range("D26").Formula = "'28.51"
$ws.Range("E26").Value2 = "  +1.68%  "
$ws.Range("D27").Formula = "'4.43"
$ws.Range("E27").Value2 = "  +1.41%  "
$ws.Range("D28").Formula = "'7.90"
$ws.Range("E28").Value2 = "  -3.34%  "
$ws.Range("D29").Formula = "'7.34"
$ws.Range("E29").Value2 = "  -1.04%  "
$ws.Range("E30").Value2 = "  -0.47%  "
$ws.Range("E31").Value2 = "  -0.33%  "
$ws.Range("E32").Value2 = "  +2.80%  "
$ws.Range("D33").Formula = "'11.43"
$ws.Range("E33").Value2 = "  +2.61%  "
$ws.Range("D34").Formula = "'40.08"
$ws.Range("E34").Value2 = "  +9.54%  "
$ws.Range("D35").Formula = "'0.0529"
$ws.Range("E35").Value2 = "  +7.29%  "
$ws.Range("E36").Value2 = "  +0.07%  "
$ws.Range("D37").Formula = "'51.97"
$ws.Range("E37").Value2 = "  +1.40%  "
$ws.Range("D38").Formula = "'3.36"
$ws.Range("E38").Value2 = "  +8.71%  "
$ws.Range("D39").Formula = "'0.997"
$ws.Range("E39").Value2 = "  -0.34%  "
$ws.Range("D40").Formula = "'3.49"
$ws.Range("E40").Value2 = "  -0.78%  "
$ws.Range("D41").Formula = "'137.66"
$ws.Range("E41").Value2 = "  +2.01%  "
$ws.Range("E42").Value2 = "  +2.33%  "
$ws.Range("D43").Formula = "'1.89"
$ws.Range("E43").Value2 = "  -1.41%  "
$ws.Range("E44").Value2 = "  -1.72%  "
$ws.Range("D45").Formula = "'16.70"
$ws.Range("E45").Value2 = "  -3.85%  "
$ws.Range("E46").Value2 = "  -1.45%  "
$ws.Range("E47").Value2 = "  +11.29%  "
$ws.Range("D48").Formula = "'22.43"
$ws.Range("E48").Value2 = "  +1.03%  "
$ws.Range("D49").Value2 = "2.174.41"
$ws.Range("E49").Value2 = "  +2.12%  "
$ws.Range("D50").Formula = "'2.46"
$ws.Range("E50").Value2 = "  +0.75%  "
$ws.Range("E51").Value2 = "  -11.96%  "
